# Refined metadata to be additional tab
#
# 1) Update the "time_taken" timestamps on the "data" sheet (column F, rows 2-45)
#    to reflect the later panel re-query time (13:38:36.xxxxxx -> 14:19:07.xxxxxx).
# 2) Add a new "metadata" worksheet (placed right after "data") describing the
#    panel query itself: data_name / data_id / data_version /
#    data_version_created / panel_query_time / panel_get_request.

$wb = $excel.ActiveWorkbook
$data = $wb.Worksheets.Item("data")

$newTimes = @(
    "2021-10-05 14:19:07.131669",
    "2021-10-05 14:19:07.131677",
    "2021-10-05 14:19:07.131680",
    "2021-10-05 14:19:07.131683",
    "2021-10-05 14:19:07.131686",
    "2021-10-05 14:19:07.131689",
    "2021-10-05 14:19:07.131692",
    "2021-10-05 14:19:07.131694",
    "2021-10-05 14:19:07.131697",
    "2021-10-05 14:19:07.131700",
    "2021-10-05 14:19:07.131703",
    "2021-10-05 14:19:07.131705",
    "2021-10-05 14:19:07.131708",
    "2021-10-05 14:19:07.131710",
    "2021-10-05 14:19:07.131713",
    "2021-10-05 14:19:07.131715",
    "2021-10-05 14:19:07.131718",
    "2021-10-05 14:19:07.131721",
    "2021-10-05 14:19:07.131724",
    "2021-10-05 14:19:07.131726",
    "2021-10-05 14:19:07.131729",
    "2021-10-05 14:19:07.131731",
    "2021-10-05 14:19:07.131734",
    "2021-10-05 14:19:07.131736",
    "2021-10-05 14:19:07.131740",
    "2021-10-05 14:19:07.131742",
    "2021-10-05 14:19:07.131745",
    "2021-10-05 14:19:07.131747",
    "2021-10-05 14:19:07.131750",
    "2021-10-05 14:19:07.131753",
    "2021-10-05 14:19:07.131755",
    "2021-10-05 14:19:07.131758",
    "2021-10-05 14:19:07.131761",
    "2021-10-05 14:19:07.131764",
    "2021-10-05 14:19:07.131766",
    "2021-10-05 14:19:07.131769",
    "2021-10-05 14:19:07.131772",
    "2021-10-05 14:19:07.131774",
    "2021-10-05 14:19:07.131777",
    "2021-10-05 14:19:07.131780",
    "2021-10-05 14:19:07.131783",
    "2021-10-05 14:19:07.131785",
    "2021-10-05 14:19:07.131788",
    "2021-10-05 14:19:07.131790"
)

for ($i = 0; $i -lt $newTimes.Length; $i++) {
    $row = $i + 2
    $data.Cells.Item($row, 6).Value = $newTimes[$i]
}

# Add the new "metadata" sheet right after "data"
$meta = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $data)
$meta.Name = "metadata"

# Header row (bold / centered / top-aligned / thin-bordered), matching the
# look of the "data" sheet's own header row - copy its formatting across.
$data.Range("B1").Copy()
$meta.Range("B1:G1").PasteSpecial(-4122)
$data.Range("A2").Copy()
$meta.Range("A2").PasteSpecial(-4122)

$meta.Range("B1").Value = "data_name"
$meta.Range("C1").Value = "data_id"
$meta.Range("D1").Value = "data_version"
$meta.Range("E1").Value = "data_version_created"
$meta.Range("F1").Value = "panel_query_time"
$meta.Range("G1").Value = "panel_get_request"

$meta.Range("A2").Value = 0
$meta.Range("B2").Value = "Albinism or congenital nystagmus"
$meta.Range("C2").Value = 511

# "1.17" must stay a literal text value (it looks numeric, and a plain
# .Value assignment would get auto-coerced to the number 1.17). Stage it in
# a scratch cell formatted as Text, then copy just the *value* over so the
# destination cell picks up no stray number-format/style.
$meta.Range("Z1").NumberFormat = "@"
$meta.Range("Z1").Value = "1.17"
$meta.Range("Z1").Copy()
$meta.Range("D2").PasteSpecial(-4163)
$meta.Range("Z1").Clear()

$meta.Range("E2").Value = "2021-01-26T11:03:14.022992Z"
$meta.Range("F2").Value = "2021-10-05 14:19:07.128258"
$meta.Range("G2").Value = "https://panelapp.genomicsengland.co.uk/api/v1/panels/511/?format=json"

# Keep "data" as the active/selected sheet (unchanged from the original
# workbook view) rather than leaving the newly-added "metadata" tab active.
$data.Activate()
$data.Range("A1").Select() | Out-Null
